$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Capture formatting donors BEFORE any destination cells are touched.
#    - B8/C8 hold "01/01/2020" as literal text (not a date) with the wrap style;
#      used to paste into B15/C15 so the engine doesn't auto-convert the date.
#    - B19/C19 already carry the column B/C wrap styles used for the brand new
#      B18/C18 cells.
# ---------------------------------------------------------------------------
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)   # xlPasteValues - keeps B15's own style, pastes text

$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4122)   # xlPasteFormats - gives new cell the column-B wrap style

$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)   # gives new A13 the column-A bold style

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Row 10 - "Objetivos:" value text is replaced by the docente line (the
#    shared string that B10/C10 point at gets new content upstream).
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

# ---------------------------------------------------------------------------
# 3. Row 13 - label "Docentes responsáveis:" value moves up/out; row becomes
#    "Programa resumido:" / "Semestral".
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# ---------------------------------------------------------------------------
# 3. Row 14 - "Short syllabus:" short text
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "1. Environmental problems, causes and solutions2. Ecosystems: what they are and how they work3. Biodiversity and evolution4. Biodiversity, species interactions and population control5. The human population and its impact"
$ws.Range("C14").Value = "1. Environmental problems, causes and solutions2. Ecosystems: what they are and how they work3. Biodiversity and evolution4. Biodiversity, species interactions and population control5. The human population and its impact"

# ---------------------------------------------------------------------------
# 4. Row 15 - "Programa:" / "01/01/2020" (text already pasted above)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Programa:"
$ws.Rows(15).RowHeight = 120

# ---------------------------------------------------------------------------
# 5. Row 16 - "Syllabus:" short text
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1. Environmental problems, causes and solutions2. Ecosystems: what they are and how they work3. Biodiversity and evolution4. Biodiversity, species interactions and population control5. The human population and its impact"
$ws.Range("C16").Value = "1. Environmental problems, causes and solutions2. Ecosystems: what they are and how they work3. Biodiversity and evolution4. Biodiversity, species interactions and population control5. The human population and its impact"

# ---------------------------------------------------------------------------
# 6. Row 17 - "Avaliação:" label only, B/C cleared, height reset to default
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Rows(17).AutoFit()

# ---------------------------------------------------------------------------
# 7. Row 18 - "Método:" / docente text (new B18/C18 cells already formatted)
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Rows(18).RowHeight = 60

# ---------------------------------------------------------------------------
# 8. Row 19 - "Critério:" label; B/C text unchanged
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "Critério:"

# ---------------------------------------------------------------------------
# 9. Row 20 - "Norma de recuperação:" label; B/C text unchanged
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "Norma de recuperação:"

# ---------------------------------------------------------------------------
# 10. Row 21 - "Bibliografia:" label; B/C text unchanged; height grows
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows(21).RowHeight = 120

# ---------------------------------------------------------------------------
# 11. Row 22 removed entirely (old Bibliografia / long bibliography text row)
# ---------------------------------------------------------------------------
$ws.Rows(22).Delete()
